# STX and STY operations
# Updates the 6502 instruction reference sheet: marks Zero Page mode for
# STX/STY as Implemented ("I") and clears the (incorrect) "N" markers in
# the Absolute / Indirect X / Indirect Y columns for those two rows, which
# are not valid addressing modes for STX/STY.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 52 = STX, Row 53 = STY
# Column D = Zero Page Mode, G = Absolute Mode, J = Indirect X Mode, K = Indirect Y Mode
$ws.Range("D52").Value = "I"
$ws.Range("G52").Value = "-"
$ws.Range("J52").Value = "-"
$ws.Range("K52").Value = "-"

$ws.Range("D53").Value = "I"
$ws.Range("G53").Value = "-"
$ws.Range("J53").Value = "-"
$ws.Range("K53").Value = "-"

# New helper cell summarizing completed vs not-implemented difference
$ws.Range("C1").Formula = "=B2-B1"

# Reset the view back to the top-left of the sheet / move the selection
$ws.Range("F1").Select()
